$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.075165666666667
$ws.Range("H2").Value = 9.225497000000001
$ws.Range("I2").Value = 0.02641273658732285
$ws.Range("J2").Value = 0.02641273658732285
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05968133333333333
$ws.Range("N2").Value = 0.179044
$ws.Range("O2").Value = 0.02602747651633847
$ws.Range("P2").Value = 0.02602747651633848
$ws.Range("Q2").Value = 0.1835299872075556
$ws.Range("R2").Value = 1.651769884868
$ws.Range("S2").Value = 0.0006874568812587793
$ws.Range("T2").Value = 0.0006874568812587794
$ws.Range("G3").Value = 3.075165666666667
$ws.Range("H3").Value = 9.225497000000001
$ws.Range("I3").Value = 0.02641273658732285
$ws.Range("J3").Value = 0.02641273658732285
$ws.Range("O3").Value = 0.144012433133819
$ws.Range("P3").Value = 0.144012433133819
$ws.Range("Q3").Value = 1.015488381833222
$ws.Range("R3").Value = 9.139395436499001
$ws.Range("S3").Value = 0.003803762461663006
$ws.Range("T3").Value = 0.003803762461663006
$ws.Range("G4").Value = 3.075165666666667
$ws.Range("H4").Value = 9.225497000000001
$ws.Range("I4").Value = 0.02641273658732285
$ws.Range("J4").Value = 0.02641273658732285
$ws.Range("O4").Value = 0.8299600903498424
$ws.Range("P4").Value = 0.8299600903498425
$ws.Range("Q4").Value = 5.852375456724333
$ws.Range("R4").Value = 52.67137911051901
$ws.Range("S4").Value = 0.02192151724440106
$ws.Range("T4").Value = 0.02192151724440106
$ws.Range("I5").Value = 0.549422396165273
$ws.Range("J5").Value = 0.5494223961652731
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05968133333333333
$ws.Range("N5").Value = 0.179044
$ws.Range("O5").Value = 0.02602747651633847
$ws.Range("P5").Value = 0.02602747651633848
$ws.Range("Q5").Value = 3.817684131531999
$ws.Range("R5").Value = 34.35915718378799
$ws.Range("S5").Value = 0.01430007851374206
$ws.Range("T5").Value = 0.01430007851374206
$ws.Range("I6").Value = 0.549422396165273
$ws.Range("J6").Value = 0.5494223961652731
$ws.Range("O6").Value = 0.144012433133819
$ws.Range("P6").Value = 0.144012433133819
$ws.Range("S6").Value = 0.07912365608997399
$ws.Range("T6").Value = 0.07912365608997401
$ws.Range("I7").Value = 0.549422396165273
$ws.Range("J7").Value = 0.5494223961652731
$ws.Range("O7").Value = 0.8299600903498424
$ws.Range("P7").Value = 0.8299600903498425
$ws.Range("S7").Value = 0.4559986615615569
$ws.Range("T7").Value = 0.4559986615615571
$ws.Range("I8").Value = 0.424164867247404
$ws.Range("J8").Value = 0.4241648672474041
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.05968133333333333
$ws.Range("N8").Value = 0.179044
$ws.Range("O8").Value = 0.02602747651633847
$ws.Range("P8").Value = 0.02602747651633848
$ws.Range("Q8").Value = 2.947327036804444
$ws.Range("R8").Value = 26.52594333124
$ws.Range("S8").Value = 0.01103994112133763
$ws.Range("T8").Value = 0.01103994112133764
$ws.Range("I9").Value = 0.424164867247404
$ws.Range("J9").Value = 0.4241648672474041
$ws.Range("O9").Value = 0.144012433133819
$ws.Range("P9").Value = 0.144012433133819
$ws.Range("S9").Value = 0.06108501458218199
$ws.Range("T9").Value = 0.061085014582182
$ws.Range("I10").Value = 0.424164867247404
$ws.Range("J10").Value = 0.4241648672474041
$ws.Range("O10").Value = 0.8299600903498424
$ws.Range("P10").Value = 0.8299600903498425
$ws.Range("S10").Value = 0.3520399115438844
$ws.Range("T10").Value = 0.3520399115438845
